$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.453.35'
$ws.Range("E2").Value = '  -1.13%  '
$ws.Range("D3").Value = '2.529.40'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.58%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.588'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.04%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.539'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0811'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.68'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.92%  '
$ws.Range("E13").Value = '  -1.71%  '
$ws.Range("D14").Value = '2.902.06'
$ws.Range("E14").Value = '  -2.08%  '
$ws.Range("D15").Value = '2.523.23'
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.56%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.864'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.68%  '
$ws.Range("D18").Value = '42.413.88'
$ws.Range("E18").Value = '  -1.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.36%  '
$ws.Range("D20").Value = '0.0₃0975'
$ws.Range("E20").Value = '  -2.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.99'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.38'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.31%  '
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.04'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.35%  '
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("E28").Value = '  +9.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.30'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.34'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0791'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.93%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.115'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.20%  '
$ws.Range("B39").Value = 'Stellar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.120'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.48%  '
$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.10'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.41'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0300'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.31%  '
$ws.Range("D46").Value = '2.038.57'
$ws.Range("E46").Value = '  -1.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '84.89'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.97'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("D49").Value = '2.763.90'
$ws.Range("E49").Value = '  -1.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.96%  '
